# "updated bold in Rmd"
# The sentence "Anyone whose funding runs out between the 2 lines should
# actively be writing and submitting grants." gets "writing and
# submitting" rendered in bold (the surrounding run is split into three
# runs, with only the middle one turning bold).

$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("writing and submitting", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Font.Bold = $true
} else {
    throw "Could not find target text 'writing and submitting'"
}
